$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values to append to column A, starting at row 106 (rows 106-165).
$newValues = @(
    '0.95 UF 25-02-21',
    '0.95 UF 25-02-21',
    '0.95 UF 02-03-21',
    '$ 34.594 CLP 02-03-21',
    '$ 34.594 CLP 02-03-21',
    '$ 34.594 CLP 02-03-21',
    '$ 34.594 CLP 02-03-21',
    '3.57 UF 03-03-21',
    '0.95 UF 03-03-21',
    '$ 34.594 CLP 04-03-21',
    '3.57 UF 04-03-21',
    '5.95 UF 04-03-21',
    '3.57 UF 04-03-21',
    ' 04-03-21',
    '23.80 UF 04-03-21',
    '2.38 UF 04-03-21',
    '7.14 UF 04-03-21',
    '4.76 UF 04-03-21',
    '13.09 UF 04-03-21',
    '23.80 UF 04-03-21',
    '41.65 UF 04-03-21',
    '22.61 UF 04-03-21',
    '22.61 UF 04-03-21',
    '41.65 UF 04-03-21',
    '23.80 UF 04-03-21',
    '23.80 UF 04-03-21',
    '$ 69.187 CLP 04-03-21',
    '$ 103.781 CLP 04-03-21',
    '$ 172.968 CLP 04-03-21',
    '$ 172.968 CLP 04-03-21',
    '23.80 UF 04-03-21',
    '$ 34.594 CLP 08-03-21',
    '$ 69.187 CLP 08-03-21',
    '3.57 UF 08-03-21',
    '$ 103.781 CLP 08-03-21',
    '5.95 UF 08-03-21',
    '$ 172.968 CLP 08-03-21',
    '23.80 UF 08-03-21',
    '0.95 UF 08-03-21',
    '2.38 UF 08-03-21',
    '7.14 UF 08-03-21',
    '4.76 UF 08-03-21',
    '13.09 UF 08-03-21',
    '23.80 UF 08-03-21',
    '41.65 UF 08-03-21',
    '22.61 UF 08-03-21',
    '2.38 UF 08-03-21',
    '$ 34.594 CLP 08-03-21',
    '2.38 UF 08-03-21',
    '$ 34.594 CLP 08-03-21',
    '$ 34.594 CLP 08-03-21',
    '41.65 UF 08-03-21',
    '$ 69.187 CLP 08-03-21',
    '$ 103.781 CLP 08-03-21',
    '$ 172.968 CLP 08-03-21',
    '$ 172.968 CLP 08-03-21',
    '41.65 UF 09-03-21',
    '$ 174.752 CLP 18-03-21',
    '23.80 UF 18-03-21',
    '23.80 UF 18-03-21'
)

$startRow = 106
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newValues[$i]
}
